# Regenerate the "K" column (column G) values in the save_data sheet.
# The commit replaces the old "Strike#" derived K values with freshly
# calculated ones (std/mean based s_vals calc), row by row for rows 2-45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 3
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    22 = 0
    23 = 1
    24 = 1
    25 = 0
    26 = 3
    27 = 1
    28 = 1
    29 = 1
    30 = 2
    31 = 2
    32 = 0
    33 = 0
    34 = 0
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 2
    43 = 0
    44 = 1
    45 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
